$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 371
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
